$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.888.63"
$ws.Range("E2").Value = "  -1.78%  "
$ws.Range("D3").Value = "1.832.39"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'310.58"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").Value = "'0.4617"
$ws.Range("E7").Value = "  -0.59%  "
$ws.Range("D8").Value = "'0.3669"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "'0.07175"
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").Value = "'0.8780"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "'0.07868"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "'19.62"
$ws.Range("E12").Value = "  -1.32%  "
$ws.Range("D13").Value = "1.839.60"
$ws.Range("E13").Value = "  -4.34%  "
$ws.Range("D14").Value = "'5.335"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "'6.391"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "'87.19"
$ws.Range("E16").Value = "  -5.40%  "
$ws.Range("D17").Value = "'1.007"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "'0.000008738"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "26.920.94"
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("D21").Value = "'14.47"
$ws.Range("E21").Value = "  -2.45%  "
$ws.Range("D22").Value = "'4.999"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("D23").Value = "'10.45"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("D24").Value = "'1.975"
$ws.Range("E24").Value = "  +3.93%  "
$ws.Range("D25").Value = "'150.73"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").Value = "'18.24"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "'1.972"
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("D28").Value = "'113.58"
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("D29").Value = "'4.939"
$ws.Range("E29").Value = "  -3.64%  "
$ws.Range("D30").Value = "'0.08820"
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").Value = "'3.129"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("D32").Value = "'0.7567"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").Value = "'4.461"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").Value = "'1.131"
$ws.Range("E34").Value = "  -2.59%  "
$ws.Range("D35").Value = "'2.588"
$ws.Range("E35").Value = "  -2.06%  "
$ws.Range("D36").Value = "'1.093"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("D37").Value = "'0.01934"
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").Value = "'0.05135"
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("D40").Value = "'6.911"
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("D41").Value = "'0.4982"
$ws.Range("E41").Value = "  -3.50%  "
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").Value = "'8.360"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "'0.4676"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "'1.006"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'10.16"
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").Value = "'102.33"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").Value = "'64.46"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").Value = "'36.41"
$ws.Range("E51").Value = "  -1.88%  "
